# Auto-generated edit script applying cached market-data value updates
# (H:N columns) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ALC!33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 220.78378
$ws.Range("I33").Value = 155.74193
$ws.Range("J33").Value = 556.8333
$ws.Range("K33").Value = 155.74193
$ws.Range("L33").Value = 556.8333
$ws.Range("M33").Value = 73.25807
$ws.Range("N33").Value = -1014.8333

# ALC!40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1963.7391
$ws.Range("I40").Value = 1640.2
$ws.Range("K40").Value = 1640.2
$ws.Range("M40").Value = -1465.2

# ALC!86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7286.8423
$ws.Range("I86").Value = 9779.091
$ws.Range("J86").Value = 3860
$ws.Range("K86").Value = 9779.091
$ws.Range("L86").Value = 3860
$ws.Range("M86").Value = -8656.091
$ws.Range("N86").Value = -6106

# ALC!89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7286.8423
$ws.Range("I89").Value = 9779.091
$ws.Range("J89").Value = 3860
$ws.Range("K89").Value = 48895.455
$ws.Range("L89").Value = 19300
$ws.Range("M89").Value = -43279.455
$ws.Range("N89").Value = -30532

# ALC!137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 38462820
$ws.Range("I137").Value = 62500896
$ws.Range("J137").Value = 1898.6
$ws.Range("K137").Value = 187502688
$ws.Range("L137").Value = 5695.799999999999
$ws.Range("M137").Value = -187500138
$ws.Range("N137").Value = -10795.8

# ARM!2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 885.06976
$ws.Range("I2").Value = 799.9375
$ws.Range("J2").Value = 1132.7273
$ws.Range("K2").Value = 799.9375
$ws.Range("L2").Value = 1132.7273
$ws.Range("M2").Value = -686.9375
$ws.Range("N2").Value = -1358.7273

# ARM!45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 78756.62
$ws.Range("I45").Value = 112215.336
$ws.Range("J45").Value = 3474.5
$ws.Range("K45").Value = 112215.336
$ws.Range("L45").Value = 3474.5
$ws.Range("M45").Value = -111838.336
$ws.Range("N45").Value = -4228.5

# ARM!74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8621836
$ws.Range("I74").Value = 10205226
$ws.Range("J74").Value = 1160.2222
$ws.Range("K74").Value = 10205226
$ws.Range("L74").Value = 1160.2222
$ws.Range("M74").Value = -10204352
$ws.Range("N74").Value = -2908.2222

# ARM!77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8621836
$ws.Range("I77").Value = 10205226
$ws.Range("J77").Value = 1160.2222
$ws.Range("K77").Value = 51026130
$ws.Range("L77").Value = 5801.111
$ws.Range("M77").Value = -51021762
$ws.Range("N77").Value = -14537.111

# ARM!116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 885.06976
$ws.Range("I116").Value = 799.9375
$ws.Range("J116").Value = 1132.7273
$ws.Range("K116").Value = 799.9375
$ws.Range("L116").Value = 1132.7273
$ws.Range("M116").Value = 1494.0625
$ws.Range("N116").Value = -5720.7273

# ARM!132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14495335
$ws.Range("I132").Value = 18868710
$ws.Range("J132").Value = 8530.9375
$ws.Range("K132").Value = 56606130
$ws.Range("L132").Value = 25592.8125
$ws.Range("M132").Value = -56603600
$ws.Range("N132").Value = -30652.8125

# BSM!3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 885.06976
$ws.Range("I3").Value = 799.9375
$ws.Range("J3").Value = 1132.7273
$ws.Range("K3").Value = 799.9375
$ws.Range("L3").Value = 1132.7273
$ws.Range("M3").Value = -685.9375
$ws.Range("N3").Value = -1360.7273

# BSM!17
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# CRP!48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 6999.6665
$ws.Range("J48").Value = 6999.6665
$ws.Range("L48").Value = 6999.6665
$ws.Range("N48").Value = -7951.6665

# CRP!58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 50000976
$ws.Range("I58").Value = 125000970
$ws.Range("J58").Value = 981.6667
$ws.Range("K58").Value = 125000970
$ws.Range("L58").Value = 981.6667
$ws.Range("M58").Value = -125000767
$ws.Range("N58").Value = -1387.6667

# CRP!94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 7193.4443
$ws.Range("I94").Value = 17248.666
$ws.Range("K94").Value = 17248.666
$ws.Range("M94").Value = -16797.666

# CRP!105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2323.8987
$ws.Range("I105").Value = 2251.4211
$ws.Range("K105").Value = 2251.4211
$ws.Range("M105").Value = -504.4211

# CRP!132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 16673462
$ws.Range("I132").Value = 1578.125
$ws.Range("J132").Value = 27788052
$ws.Range("K132").Value = 4734.375
$ws.Range("L132").Value = 83364156
$ws.Range("M132").Value = -2204.375
$ws.Range("N132").Value = -83369216

# CRP!134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1007
$ws.Range("I134").Value = 1026.6316
$ws.Range("J134").Value = 965.55554
$ws.Range("K134").Value = 3079.8948
$ws.Range("L134").Value = 2896.66662
$ws.Range("M134").Value = -544.8948
$ws.Range("N134").Value = -7966.66662

# CRP!136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 50000976
$ws.Range("I136").Value = 125000970
$ws.Range("J136").Value = 981.6667
$ws.Range("K136").Value = 375002910
$ws.Range("L136").Value = 2945.0001
$ws.Range("M136").Value = -375000360
$ws.Range("N136").Value = -8045.0001

# GSM!14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 111359.89
$ws.Range("I14").Value = 125167.375
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 125167.375
$ws.Range("L14").Value = 900
$ws.Range("M14").Value = -124999.375
$ws.Range("N14").Value = -1236

# GSM!24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3000
$ws.Range("J24").Value = 3000
$ws.Range("L24").Value = 3000
$ws.Range("N24").Value = -3346

# GSM!26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 6133.3335
$ws.Range("J26").Value = 6133.3335
$ws.Range("L26").Value = 6133.3335
$ws.Range("N26").Value = -6693.3335

# GSM!50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 6133.3335
$ws.Range("J50").Value = 6133.3335
$ws.Range("L50").Value = 6133.3335
$ws.Range("N50").Value = -7129.3335

# GSM!141
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 40302.8
$ws.Range("J141").Value = 40302.8
$ws.Range("L141").Value = 40302.8
$ws.Range("N141").Value = -50662.8

# LTW!122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5671.9395
$ws.Range("I122").Value = 6516
$ws.Range("J122").Value = 3421.111
$ws.Range("K122").Value = 19548
$ws.Range("L122").Value = 10263.333
$ws.Range("M122").Value = -17098
$ws.Range("N122").Value = -15163.333

# LTW!123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 29538.611
$ws.Range("J123").Value = 29538.611
$ws.Range("L123").Value = 29538.611
$ws.Range("N123").Value = -39338.611

# LTW!132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18874348
$ws.Range("I132").Value = 40001744
$ws.Range("J132").Value = 10599.75
$ws.Range("K132").Value = 120005232
$ws.Range("L132").Value = 31799.25
$ws.Range("M132").Value = -120002702
$ws.Range("N132").Value = -36859.25

# WVR!82
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 38950
$ws.Range("I82").Value = 17800
$ws.Range("J82").Value = 46000
$ws.Range("K82").Value = 17800
$ws.Range("L82").Value = 46000
$ws.Range("M82").Value = -17417
$ws.Range("N82").Value = -46766

# WVR!85
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 38950
$ws.Range("I85").Value = 17800
$ws.Range("J85").Value = 46000
$ws.Range("K85").Value = 17800
$ws.Range("L85").Value = 46000
$ws.Range("M85").Value = -16474
$ws.Range("N85").Value = -48652

# WVR!135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 36852.555
$ws.Range("J135").Value = 36852.555
$ws.Range("L135").Value = 36852.555
$ws.Range("N135").Value = -46992.555

# WVR!136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3018.982
$ws.Range("I136").Value = 3699.9707
$ws.Range("J136").Value = 1916.4286
$ws.Range("K136").Value = 11099.9121
$ws.Range("L136").Value = 5749.2858
$ws.Range("M136").Value = -8549.9121
$ws.Range("N136").Value = -10849.2858

# WVR!141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 45485.8
$ws.Range("J141").Value = 45485.8
$ws.Range("L141").Value = 45485.8
$ws.Range("N141").Value = -55845.8
